$d = $word.ActiveDocument

# Set the header row (row 1) cells of every table to bottom vertical alignment.
# wdCellAlignVerticalBottom = 3
foreach ($t in $d.Tables) {
    $row = $t.Rows.Item(1)
    foreach ($cell in $row.Cells) {
        $cell.VerticalAlignment = 3
    }
}
